# Add "unsold" (未售) / "total" (总量) columns to every floor-plan sheet.
#
# Each worksheet (except the empty "First" cover sheet) already has a
# "去化率" (sell-through rate) header in F1 and a computed rate in F2.
# This change appends, per-sheet, two new header cells in row 4
# (F4 = "未售", G4 = "总量") and two new numeric value cells in row 5
# (F5 = unsold-unit count, G5 = total-unit count) using the same visual
# style (green fill, centered + wrapped) already used by F1/F2.

$wb = $excel.ActiveWorkbook

# Per-sheet (F5, G5) = (unsold units, total units) values.
$values = @{
    "1"  = @(10, 72)
    "2"  = @(29, 58)
    "3"  = @(10, 68)
    "4"  = @(9, 72)
    "5"  = @(17, 72)
    "6"  = @(5, 68)
    "7"  = @(13, 59)
    "8"  = @(10, 72)
    "9"  = @(89, 108)
    "10" = @(106, 106)
    "11" = @(106, 106)
    "12" = @(107, 108)
}

foreach ($ws in $wb.Worksheets) {
    $name = [string]$ws.Name
    if (-not $values.ContainsKey($name)) {
        continue
    }

    $pair = $values[$name]
    $unsold = $pair[0]
    $total = $pair[1]

    $ws.Cells.Item(4, 6).Value = "未售"
    $ws.Cells.Item(4, 7).Value = "总量"
    $ws.Cells.Item(5, 6).Value = $unsold
    $ws.Cells.Item(5, 7).Value = $total

    # Match the existing F1/F2 styling (green fill, centered + wrapped)
    # by copying the format from F1 onto the new F4:G5 block.
    $ws.Cells.Item(1, 6).Copy()
    $fmtRange = $ws.Range("F4:G5")
    $fmtRange.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0

Write-Output "done"
